# A new weekly price record was inserted into the "Pomelo" price log.
# It belongs chronologically right after the existing row 622, so it is
# inserted as a new row 623 - pushing the former rows 623:687 down to
# 624:688 - and then populated with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 623, shifting rows 623:687 down to 624:688.
$ws.Rows("623:623").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A623").Value = 4
$ws.Range("B623").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C623").Value = "Los Lagos"
$ws.Range("D623").Value = 45212
$ws.Range("E623").Value = 10
$ws.Range("F623").Value = "Fruta"
$ws.Range("G623").Value = 100102
$ws.Range("H623").Value = "Cítricos"
$ws.Range("I623").Value = 100102006
$ws.Range("J623").Value = "Pomelo"
$ws.Range("K623").Value = "Start Ruby"
$ws.Range("L623").Value = "Primera"
$ws.Range("M623").Value = 150
$ws.Range("N623").Value = 15000
$ws.Range("O623").Value = 15000
$ws.Range("P623").Value = 15000
$ws.Range("Q623").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R623").Value = "Región de O'Higgins"
$ws.Range("S623").Value = 1071
$ws.Range("T623").Value = 14
